$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "2022-04-13"
$ws.Range("I2").Style = "Normal"

$ws.Range("J2").Value = "NIFTY2241318000CE"
$ws.Range("AA2").Value = "existing"

# Row 3 updates
$ws.Range("A3").Value = "Sell"

$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "2022-04-13"
$ws.Range("I3").Style = "Normal"

$ws.Range("J3").Value = "NIFTY2241318000CE"
$ws.Range("AA3").Value = "existing"
